# Update mapa interactivo PEBCOM: remove obsolete/duplicate rows.
# Rows being removed (identified by their unique "Caso" / OT values before the edit):
#   Row 11 -> Caso 2800, Anchorena 1288
#   Row 67 -> Caso 6137, LA PLATA AV. 1058
#   Row 75 -> Caso 6522, ESTADO PLURINACIONAL DE BOLIVIA 384
#
# Deleting from the bottom-most row first keeps the remaining row numbers stable
# while the deletions are performed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(75).Delete()
$ws.Rows.Item(67).Delete()
$ws.Rows.Item(11).Delete()
